# Apply the "updates thurs from @mac" revision:
#  - Tab Banners: move the saved cell selection from B6 to A7 (no data change)
#  - Terminology: append a new glossary row (29) describing the new
#    "Pivot Table Report" download option, copying the look of the
#    existing "Combined Report (Master)" row (A-col shading) and the
#    wrapped-description look used elsewhere in column B, then leave the
#    sheet's saved selection on B32:B33 like the source file shows.

$wb = $excel.ActiveWorkbook

# --- Sheet "Tab Banners" ---------------------------------------------------
$wsBanners = $wb.Worksheets.Item("Tab Banners")
$wsBanners.Range("A7").Select()

# --- Sheet "Terminology" ----------------------------------------------------
$wsTerm = $wb.Worksheets.Item("Terminology")
$wsTerm.Activate()

$wsTerm.Range("A29").Value = "Pivot Table Report"
$wsTerm.Range("B29").Value = "The output of this report includes summary tables (pivot tables) for GFDRR's Porfolio by: `n(1) GP/Global Theme & Disbursment Risk Level; `n(2) Trustee & Disbursement Risk Level;`n(3) Country/Region and Disbursement Risk Level. `nThe report can be customized by selecting (using the buttons on the right) relevant categories to include (or not) in the report.  The report also provides a raw list of the grants used for the tables, for your reference."

# Match the shading already used for the other "report" rows (25-28) and the
# wrapped-text look used for long descriptions (e.g. B18).
$wsTerm.Range("A28").Copy()
$wsTerm.Range("A29").PasteSpecial(-4122)
$wsTerm.Range("B18").Copy()
$wsTerm.Range("B29").PasteSpecial(-4122)
$wsTerm.Range("A1").Select()

$wsTerm.Rows.Item(29).RowHeight = 119

$wsTerm.Range("B32:B33").Select()
